$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 443, pushing existing rows 443:493 down to 444:494.
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new record.
$ws.Range("A443").Value2 = 3
$ws.Range("B443").Value2 = "Femacal de La Calera"
$ws.Range("C443").Value2 = "Coquimbo"
$ws.Range("D443").Value2 = 44946
$ws.Range("E443").Value2 = 5
$ws.Range("F443").Value2 = 100112040
$ws.Range("G443").Value2 = "Cilantro"
$ws.Range("H443").Value2 = "Sin especificar"
$ws.Range("I443").Value2 = "Primera"
$ws.Range("J443").Value2 = 270
$ws.Range("K443").Value2 = 4000
$ws.Range("L443").Value2 = 4500
$ws.Range("M443").Value2 = 4241
$ws.Range("N443").Value2 = "$/docena de atados (3 kilos)"
$ws.Range("O443").Value2 = "Provincia de Quillota"
$ws.Range("P443").Value2 = 1414
$ws.Range("Q443").Value2 = 3
$ws.Range("R443").Value2 = "Hortaliza"
